# Updated cryptos list on Mon Apr  8 06:36:23 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'69.753.01"
$ws.Range("E2").Value = "  +0.80%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'3.426.42"
$ws.Range("E3").Value = "  +1.21%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'583.18"
$ws.Range("E5").Value = "  -0.69%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'176.74"
$ws.Range("E6").Value = "  -1.88%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "'3.419.98"
$ws.Range("E7").Value = "  +1.19%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.06%  "

# Row 9 - XRP
$ws.Range("D9").Value = "'0.595"
$ws.Range("E9").Value = "  +0.04%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +2.82%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  -1.08%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "'49.02"
$ws.Range("E12").Value = "  +0.85%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +0.45%  "

# Row 14 - BitcoinCash
$ws.Range("D14").Value = "'690.68"
$ws.Range("E14").Value = "  +2.09%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "'3.975.30"
$ws.Range("E15").Value = "  +1.10%  "

# Row 16 - Polkadot
$ws.Range("E16").Value = "  +0.11%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "'69.778.58"
$ws.Range("E17").Value = "  +0.76%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "'3.426.38"
$ws.Range("E18").Value = "  +1.00%  "

# Row 19 - TRON
$ws.Range("E19").Value = "  +1.14%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  +0.11%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "'11.41"
$ws.Range("E21").Value = "  +0.73%  "

# Row 22 - Polygon
$ws.Range("D22").Value = "'0.897"
$ws.Range("E22").Value = "  -0.36%  "

# Row 23 - Toncoin
$ws.Range("D23").Value = "'5.48"
$ws.Range("E23").Value = "  +1.20%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").Value = "'16.94"
$ws.Range("E24").Value = "  -0.93%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "'100.83"
$ws.Range("E25").Value = "  -2.63%  "

# Row 26 - PancakeSwap
$ws.Range("E26").Value = "  +0.19%  "

# Row 27 - ImmutableX
$ws.Range("D27").Value = "'2.65"
$ws.Range("E27").Value = "  -2.40%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "'9.60"
$ws.Range("E28").Value = "  +0.12%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "'33.46"
$ws.Range("E29").Value = "  -1.96%  "

# Row 30 - Filecoin
$ws.Range("D30").Value = "'8.76"
$ws.Range("E30").Value = "  +0.64%  "

# Row 31 - NEARProtocol
$ws.Range("D31").Value = "'7.15"
$ws.Range("E31").Value = "  +2.34%  "

# Row 32 - Bittensor
$ws.Range("D32").Value = "'573.98"
$ws.Range("E32").Value = "  +3.69%  "

# Row 33 - dogwifhat
$ws.Range("E33").Value = "  +2.43%  "

# Row 34 - Cosmos
$ws.Range("D34").Value = "'11.02"
$ws.Range("E34").Value = "  -1.47%  "

# Row 35 - OKB
$ws.Range("D35").Value = "'58.44"
$ws.Range("E35").Value = "  +0.85%  "

# Row 37 - Dai
$ws.Range("E37").Value = "  +0.06%  "

# Row 38 - Maker
$ws.Range("D38").Value = "'3.608.24"
$ws.Range("E38").Value = "  -2.23%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  -0.22%  "

# Row 40 - InjectiveProtocol
$ws.Range("D40").Value = "'35.15"
$ws.Range("E40").Value = "  +0.26%  "

# Row 41 - PEPE
$ws.Range("D41").Value = "'0.0₃0739"
$ws.Range("E41").Value = "  +5.28%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  +0.56%  "

# Row 43 - Fetch.AI
$ws.Range("D43").Value = "'2.67"
$ws.Range("E43").Value = "  +0.34%  "

# Row 44 - was ApeXProtocol, now VeChain
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0421"
$ws.Range("E44").Value = "  -0.43%  "

# Row 45 - TheGraph
$ws.Range("E45").Value = "  -1.34%  "

# Row 46 - was VeChain, now Mantle
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'1.46"
$ws.Range("E46").Value = "  +4.35%  "

# Row 47 - was Mantle, now ThetaToken
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Value = "'2.67"
$ws.Range("E47").Value = "  +0.32%  "

# Row 48 - was ThetaToken, now Stellar
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.129"
$ws.Range("E48").Value = "  -0.49%  "

# Row 49 - was Stellar, now FirstDigitalUSD
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").Value = "'0.999"
$ws.Range("E49").Value = "  -0.20%  "

# Row 50 - was FirstDigitalUSD, now Monero
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'133.15"
$ws.Range("E50").Value = "  +0.99%  "

# Row 51 - was Monero, now CoreDAO
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").Value = "'2.66"
$ws.Range("E51").Value = "  +2.32%  "
